$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("A3").Value = 130937854
$ws.Range("B3").Value = 57881
$ws.Range("E3").Value = 100049
$ws.Range("F3").Value = "Spillkråka"
$ws.Range("G3").Value = "Dryocopus martius"
$ws.Range("Q3").Value = 489668
$ws.Range("R3").Value = 7004128
$ws.Range("AC3").Value = "Rejäla hackspår, färska och äldre, I två levande granar och i ytlig grov rotdel."

# --- Row 4 ---
$ws.Range("A4").Value = 130937843
$ws.Range("M4").Value = "färska spår"
$ws.Range("Q4").Value = 489760
$ws.Range("R4").Value = 7004232
$ws.Range("AC4").Value = "Ringhack, färska och äldre, i riklig mängd längs flera meter högt upp på en granstam med spår av rikligt sav/kådaflöde."

# --- Row 5 ---
$ws.Range("A5").Value = 130937852
$ws.Range("B5").Value = 57884
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("M5").Value = "äldre spår"
$ws.Range("Q5").Value = 489520
$ws.Range("R5").Value = 7004161
$ws.Range("AC5").Value = "Ringhack, äldre, ytliga enstaka längs flera meter på en granstam vid kanten mot yngre skog."

# --- Row 19 ---
$ws.Range("J19").ClearContents()
$ws.Range("AF19").ClearContents()
$ws.Range("A19").Value = 130937846
$ws.Range("B19").Value = 57884
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 100109
$ws.Range("F19").Value = "Tretåig hackspett"
$ws.Range("G19").Value = "Picoides tridactylus"
$ws.Range("H19").Value = "(Linnaeus, 1758)"
$ws.Range("M19").Value = "färska spår"
$ws.Range("Q19").Value = 489591
$ws.Range("R19").Value = 7004206
$ws.Range("AC19").Value = "Ringhack, främst färska och några äldre, längs flera meter på en granstam vid kant mot yngre skog. Fyndplatsen finns i barrblandskog intill några rotvältor."
$ws.Range("AJ19").Value = "gran"
$ws.Range("AK19").Value = "Picea abies"
$ws.Range("AM19").Value = "Trädstam på levande träd"
$ws.Range("AO19").Value = "Stem on living tree # Picea abies"

# --- Row 20 ---
$ws.Range("M20").ClearContents()
$ws.Range("AJ20").ClearContents()
$ws.Range("AK20").ClearContents()
$ws.Range("AM20").ClearContents()
$ws.Range("AO20").ClearContents()
$ws.Range("J20").Value = ""
$ws.Range("AF20").Value = ""
$ws.Range("A20").Value = 130937860
$ws.Range("B20").Value = 97879
$ws.Range("D20").Value = "LC"
$ws.Range("E20").Value = 221945
$ws.Range("F20").Value = "Revlummer"
$ws.Range("G20").Value = "Lycopodium annotinum"
$ws.Range("H20").Value = "L."
$ws.Range("Q20").Value = 489614
$ws.Range("R20").Value = 7004216
$ws.Range("AC20").Value = "Växer här i barrblandskog."
